$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3-17 (no longer part of the data set) by deleting the entire rows,
# which shifts everything below up and shrinks the used range accordingly.
$ws.Range("A3:B17").EntireRow.Delete()

# Update the remaining data row (row 2) with the new values.
# The Number column holds a long numeric-looking string (not an actual number),
# so prefix with an apostrophe to force it to be stored as text, matching the
# original cell's text formatting.
$ws.Range("A2").Value = "'9358883639"
$ws.Range("B2").Value = "Mani agah"
